$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.401845812797546
$ws.Range("B1").Value = 1.841027736663818
$ws.Range("C1").Value = 2.895308494567871
$ws.Range("D1").Value = 4.753347396850586
$ws.Range("E1").Value = 1.044169306755066
